$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A ("Confirmation Number" header, but containing junk
# values "1"/"1"/":1") is removed entirely. This shifts CRS Number -> A,
# Name -> B, Booking.com Price -> C, Description -> D.
$ws.Columns.Item(1).Delete()

# Re-label the headers to match the new sheet.
$ws.Range("A1").Value = "Conf number"
$ws.Range("B1").Value = "Guest Name"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Description"

# Update the description text for the first two data rows.
$ws.Range("D2").Value = "Date changed"
$ws.Range("D3").Value = "Date changed"
